$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.752.46'
$ws.Range("E2").Value = '  +0.69%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.121.04'
$ws.Range("E3").Value = '  +10.56%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '256.41'
$ws.Range("E5").Value = '  +2.29%  '

# Row 6
$ws.Range("E6").Value = '  -3.48%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.24'
$ws.Range("E8").Value = '  +6.46%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.30'
$ws.Range("E9").Value = '  +4.91%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.375'
$ws.Range("E10").Value = '  +2.50%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0742'
$ws.Range("E11").Value = '  -3.08%  '

# Row 12
$ws.Range("E12").Value = '  +0.58%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.423.59'
$ws.Range("E13").Value = '  +10.31%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.47'
$ws.Range("E14").Value = '  +0.22%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.848'
$ws.Range("E15").Value = '  +5.85%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.122.22'
$ws.Range("E16").Value = '  +10.62%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.14'
$ws.Range("E17").Value = '  +0.38%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.785.63'
$ws.Range("E18").Value = '  +0.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.90'
$ws.Range("E19").Value = '  -0.47%  '

# Row 20
$ws.Range("E20").Value = '  -1.70%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.42'
$ws.Range("E21").Value = '  +1.27%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '242.42'
$ws.Range("E22").Value = '  -3.25%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.25'
$ws.Range("E23").Value = '  +0.44%  '

# Row 24
$ws.Range("E24").Value = '  +0.09%  '

# Row 25
$ws.Range("E25").Value = '  -7.48%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.54'
$ws.Range("E26").Value = '  +2.30%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.62'
$ws.Range("E27").Value = '  +14.91%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.22'
$ws.Range("E28").Value = '  +4.90%  '

# Row 29
$ws.Range("E29").Value = '  -7.94%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.54'
$ws.Range("E30").Value = '  +54.95%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.124'
$ws.Range("E31").Value = '  -4.31%  '

# Row 32
$ws.Range("E32").Value = '  -0.16%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0957'
$ws.Range("E33").Value = '  +8.69%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0602'
$ws.Range("E34").Value = '  -1.24%  '

# Row 35
$ws.Range("E35").Value = '  +17.01%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.943'
$ws.Range("E36").Value = '  +9.21%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.89'
$ws.Range("E37").Value = '  -4.51%  '

# Row 38
$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.18%  '

# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.19'
$ws.Range("E39").Value = '  -3.45%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.37'
$ws.Range("E40").Value = '  -7.53%  '

# Row 41
$ws.Range("E41").Value = '  +7.79%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0227'
$ws.Range("E42").Value = '  -0.82%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.42'
$ws.Range("E43").Value = '  -7.68%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'
$ws.Range("E44").Value = '  +14.87%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.31'
$ws.Range("E45").Value = '  -5.17%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.360.37'
$ws.Range("E46").Value = '  +0.92%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0844'
$ws.Range("E47").Value = '  +4.34%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.14'
$ws.Range("E48").Value = '  +10.59%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.315.40'
$ws.Range("E49").Value = '  +10.02%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  -2.32%  '

# Row 51
$ws.Range("E51").Value = '  +0.75%  '
